{"js": "// Replace each two-digit-by-two-digit multiplication problem text with its\n// new value, as described by the diff. Every old value is unique within the\n// document, so a simple exact-match search/replace for each pair is safe.\nconst replacements = [\n  [\"97\u00d713=\", \"70\u00d789=\"],\n  [\"45\u00d795=\", \"29\u00d740=\"],\n  [\"19\u00d764=\", \"97\u00d765=\"],\n  [\"81\u00d756=\", \"84\u00d713=\"],\n  [\"59\u00d748=\", \"98\u00d751=\"],\n  [\"25\u00d772=\", \"78\u00d749=\"],\n  [\"53\u00d759=\", \"69\u00d721=\"],\n  [\"75\u00d726=\", \"75\u00d762=\"],\n  [\"42\u00d712=\", \"90\u00d780=\"],\n  [\"98\u00d735=\", \"57\u00d729=\"],\n  [\"50\u00d717=\", \"45\u00d770=\"],\n  [\"84\u00d743=\", \"35\u00d766=\"],\n  [\"69\u00d713=\", \"32\u00d759=\"],\n  [\"24\u00d718=\", \"93\u00d732=\"],\n  [\"20\u00d792=\", \"25\u00d717=\"],\n  [\"61\u00d738=\", \"34\u00d760=\"],\n  [\"82\u00d779=\", \"70\u00d794=\"],\n  [\"78\u00d743=\", \"12\u00d745=\"],\n  [\"78\u00d734=\", \"93\u00d795=\"],\n  [\"26\u00d791=\", \"47\u00d745=\"],\n  [\"85\u00d754=\", \"57\u00d776=\"],\n  [\"31\u00d727=\", \"61\u00d785=\"],\n  [\"52\u00d742=\", \"29\u00d713=\"],\n  [\"98\u00d767=\", \"14\u00d740=\"],\n  [\"15\u00d747=\", \"95\u00d753=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication problem text with its\n# new value, as described by the diff. Every old value is unique within the\n# document, so a Find/Replace-All pass for each pair touches exactly the\n# one cell that should change.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"97\u00d713=\", \"70\u00d789=\"),\n    @(\"45\u00d795=\", \"29\u00d740=\"),\n    @(\"19\u00d764=\", \"97\u00d765=\"),\n    @(\"81\u00d756=\", \"84\u00d713=\"),\n    @(\"59\u00d748=\", \"98\u00d751=\"),\n    @(\"25\u00d772=\", \"78\u00d749=\"),\n    @(\"53\u00d759=\", \"69\u00d721=\"),\n    @(\"75\u00d726=\", \"75\u00d762=\"),\n    @(\"42\u00d712=\", \"90\u00d780=\"),\n    @(\"98\u00d735=\", \"57\u00d729=\"),\n    @(\"50\u00d717=\", \"45\u00d770=\"),\n    @(\"84\u00d743=\", \"35\u00d766=\"),\n    @(\"69\u00d713=\", \"32\u00d759=\"),\n    @(\"24\u00d718=\", \"93\u00d732=\"),\n    @(\"20\u00d792=\", \"25\u00d717=\"),\n    @(\"61\u00d738=\", \"34\u00d760=\"),\n    @(\"82\u00d779=\", \"70\u00d794=\"),\n    @(\"78\u00d743=\", \"12\u00d745=\"),\n    @(\"78\u00d734=\", \"93\u00d795=\"),\n    @(\"26\u00d791=\", \"47\u00d745=\"),\n    @(\"85\u00d754=\", \"57\u00d776=\"),\n    @(\"31\u00d727=\", \"61\u00d785=\"),\n    @(\"52\u00d742=\", \"29\u00d713=\"),\n    @(\"98\u00d767=\", \"14\u00d740=\"),\n    @(\"15\u00d747=\", \"95\u00d753=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n"}
